$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their existing text formatting so values
# like "1.100" or "0.00001129" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.942.54"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.810.85"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "310.29"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4975"
$ws.Range("E7").Value = "  -2.75%  "
$ws.Range("D8").Value = "0.3921"
$ws.Range("E8").Value = "  +2.97%  "
$ws.Range("D9").Value = "0.09677"
$ws.Range("E9").Value = "  +24.40%  "
$ws.Range("D10").Value = "1.100"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "40.85"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "6.428"
$ws.Range("E12").Value = "  +3.84%  "
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "20.44"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "1.810.15"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "7.282"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "0.00001129"
$ws.Range("E17").Value = "  +5.17%  "
$ws.Range("D18").Value = "92.19"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "0.06646"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "17.13"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").Value = "5.913"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "28.003.84"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").Value = "2.244"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").Value = "158.82"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").Value = "2.019.29"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").Value = "20.54"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").Value = "2.385"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("D30").Value = "127.57"
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "1.034"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").Value = "5.553"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").Value = "3.613"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "0.06720"
$ws.Range("E35").Value = "  -5.07%  "
$ws.Range("D36").Value = "8.899"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").Value = "0.02325"
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("D38").Value = "0.2134"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").Value = "4.947"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "11.22"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").Value = "0.6164"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "1.156"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").Value = "13.10"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "1.289"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").Value = "0.5879"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "3.689"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").Value = "123.81"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("D49").Value = "1.930"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("D50").Value = "1.179"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").Value = "0.06766"
$ws.Range("E51").Value = "  -1.12%  "
